# Mise à jour de l'application
# Adds 7 new GPS-tracking rows (300-306) for the 2025-08-14 ("Entrainement" / "Global")
# session, and switches the numeric columns (G:V) of row 284 to right-aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 284 (G:V) -> right horizontal alignment (new style, no explicit
#    vertical alignment so it reuses the plain "horizontal=right" xf).
# ---------------------------------------------------------------------------
$ws.Range("G284:V284").VerticalAlignment = -4107   # xlBottom (keeps default, avoids adding vertical="center")
$ws.Range("G284:V284").HorizontalAlignment = -4152 # xlRight

# ---------------------------------------------------------------------------
# 2) Append the new rows 300-306.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=300; A="Entrainement"; B=45883; C="Global"; E="Emmanuel Valey";  F="left forward";     G="01:20:42"; H=4.68; I=0.3;  J=4.37; K=0.24; L=0.07; M=0;    N=0; O=0; P=3.39; Q=24.91; R=4.66; S=24; T=8; U=20; V=6 },
    @{ Row=301; A="Entrainement"; B=45883; C="Global"; E="Yoan Zouma";      F="center back";      G="01:20:09"; H=4.14; I=0.19; J=3.94; K=0.17; L=0.03; M=0;    N=0; O=0; P=2.91; Q=24.22; R=4.32; S=24; T=2; U=19; V=1 },
    @{ Row=302; A="Entrainement"; B=45883; C="Global"; E="Romain Thunet";   F="center back";      G="01:21:17"; H=5.05; I=0.09; J=4.95; K=0.09; L=0;    M=0;    N=0; O=0; P=3.66; Q=21.27; R=4.29; S=9;  T=3; U=11; V=3 },
    @{ Row=303; A="Entrainement"; B=45883; C="Global"; E="Mattheo Haon";    F="right back";       G="01:21:26"; H=5.49; I=0.35; J=5.14; K=0.29; L=0.05; M=0.01; N=0; O=1; P=3.99; Q=26.64; R=4.58; S=18; T=6; U=19; V=5 },
    @{ Row=304; A="Entrainement"; B=45883; C="Global"; E="Ilan Ihaddadene"; F="center midfield";  G="01:21:08"; H=6.4;  I=0.47; J=5.93; K=0.4;  L=0.07; M=0;    N=0; O=0; P=4.66; Q=23.98; R=4.4;  S=31; T=4; U=24; V=8 },
    @{ Row=305; A="Entrainement"; B=45883; C="Global"; E="Karahali Souaré"; F="right forward";    G="01:17:01"; H=4.7;  I=0.4;  J=4.28; K=0.24; L=0.12; M=0.06; N=0; O=6; P=3.58; Q=26.95; R=5.44; S=23; T=7; U=21; V=4 },
    @{ Row=306; A="Entrainement"; B=45883; C="Global"; E="Amine Taiar";     F="center back";      G="01:20:17"; H=4.78; I=0.19; J=4.59; K=0.16; L=0.03; M=0.01; N=0; O=1; P=3.49; Q=26.48; R=4.27; S=21; T=2; U=19; V=4 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "m/d/yy"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}

# ---------------------------------------------------------------------------
# 3) Update the active selection to match the new cursor position.
# ---------------------------------------------------------------------------
$ws.Range("D312").Select()
